# The deck's single slide master ("Integral" theme, ppt/theme/theme1.xml)
# is recoloured to the stock "Office Theme" palette (the palette that is
# already present in ppt/theme/theme2.xml, used by the notes master).
#
# The two themes only ever differed in their <a:clrScheme> (the font
# scheme and format scheme are byte-for-byte identical already), so the
# visible effect of the authored edit is reproduced by rewriting the 12
# theme colour slots exposed through Slide.ThemeColorScheme -- the
# PowerPoint object-model path that targets the presentation's live
# DrawingML theme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# msoThemeColorDark1 .. msoThemeColorFollowedHyperlink, in COM index order
# (1=dk1, 2=lt1, 3=dk2, 4=lt2, 5-10=accent1-6, 11=hlink, 12=folHlink),
# each value an RGB() style BGR-packed long, matching the "Office" theme.
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
